# Update bus voltage magnitudes (vm_pu) for the 380 kV case: rows 2-25, columns B..N (minus H).
# First column in each row is the new slack/reference voltage setpoint (1.05 -> 1.02 p.u.);
# the rest are the recomputed load-flow results for that new setpoint.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02)
$colC = @(1.03813760682279, 1.038973020341599, 1.039514084641176, 1.039741665229652, 1.039779883836698, 1.039517125125739, 1.038419835116659, 1.036490136885653, 1.03520636403418, 1.034651134305342, 1.034444996826167, 1.034489209484817, 1.034634092876963, 1.034723373556194, 1.035243226679956, 1.035569492565974, 1.035759860779142, 1.035824782032232, 1.035534480812937, 1.034591425568193, 1.033999066120522, 1.034313031877853, 1.035550300923078, 1.036988542335745)
$colD = @(1.04602409703005, 1.046696062996965, 1.047131318117771, 1.047314404946796, 1.047345152171856, 1.047133764120267, 1.046251097233499, 1.044699234376858, 1.043667121159328, 1.043220810381921, 1.043055122739619, 1.043090659107681, 1.043207112701763, 1.043278875854697, 1.043696754148765, 1.043959040318641, 1.044112085279967, 1.044164279420229, 1.043930893495849, 1.043172817465691, 1.042696718429649, 1.042949056387624, 1.043943611654897, 1.045100001694747)
$colE = @(0.992614727750844, 0.9936372048519299, 0.9942998659930998, 0.994578699834602, 0.994625531979634, 0.994303590798249, 0.9929600610674297, 0.9906006454969559, 0.989033133672735, 0.988355674866747, 0.9881042295826724, 0.9881581567098651, 0.9883348863814464, 0.9884438009545853, 0.9890781214508737, 0.989476357848556, 0.9897087662937556, 0.9897880325774034, 0.9894336180360679, 0.9882828385668249, 0.9875604150241495, 0.9879432794643023, 0.9894529299347244, 0.9912096547607049)
$colF = @(1.053209174464579, 1.054056963333636, 1.054606453087061, 1.054837675447526, 1.054876511362304, 1.054609541841754, 1.05349549838039, 1.051539491520425, 1.05024035829343, 1.049678998099634, 1.049470662059516, 1.049515342756207, 1.049661773332828, 1.049752017704458, 1.050277638797226, 1.050607662419755, 1.050800272672354, 1.050865966909931, 1.050572242300697, 1.049618648245247, 1.04902011735797, 1.049337311496263, 1.050588246778766, 1.052044315706372)
$colI = @(1.039875682886958, 1.040064517245126, 1.040185771225529, 1.040236522004472, 1.040245030105271, 1.040186450242483, 1.039939693423301, 1.039497749470177, 1.039198368914179, 1.039067615398346, 1.039018880157711, 1.039029341604098, 1.039063590342517, 1.039084669936894, 1.039207023059343, 1.039283472620222, 1.039327956270203, 1.039343105691816, 1.039275281480818, 1.039053509562866, 1.038913103619848, 1.038987627106133, 1.039278983038175, 1.039612842986546)
$colJ = @(1.043236676341394, 1.043717330985942, 1.04402817256683, 1.044158807480756, 1.044180739107755, 1.044029918286966, 1.04339915115698, 1.042286378567435, 1.041543740946479, 1.041221998512949, 1.041102463891326, 1.041128105592502, 1.041212118250062, 1.041263877923674, 1.041565090338376, 1.041753986805597, 1.041864149910714, 1.041901709761213, 1.041733721752937, 1.041187379301151, 1.040843727166638, 1.041025917003404, 1.041742878710632, 1.042574201304129)
$colK = @(1.04879069325474, 1.049274581296542, 1.049587454989105, 1.049718929730291, 1.049741001516936, 1.049589211988024, 1.048954273209874, 1.047833692852569, 1.047085547830724, 1.046761349995786, 1.046640892798331, 1.046666732870308, 1.046751393686816, 1.046803551319718, 1.047107058698096, 1.047297375878033, 1.047408360837407, 1.04744619975899, 1.047276959101509, 1.046726464161256, 1.046380140308796, 1.046563752216377, 1.047286184637782, 1.048123587710931)
$colL = @(0.9955398523335997, 0.9963617723202687, 0.9968940712668347, 0.9971179600053012, 0.9971555583673455, 0.9968970624462089, 0.9958175282591056, 0.9939188001724441, 0.9926553831429383, 0.9921088820399291, 0.9919059725120875, 0.9919494934313052, 0.9920921077337197, 0.9921799884222134, 0.9926916645766087, 0.9930127773699352, 0.9932001317071769, 0.9932640239640975, 0.9929783193494215, 0.9920501090198102, 0.9914670000341481, 0.991776070289318, 0.9929938892766442, 0.9944092447426414)
$colM = @(1.055955775347891, 1.056616474188697, 1.057044257060522, 1.057224159219025, 1.057254369187096, 1.057046660679382, 1.056179005480303, 1.05465219704054, 1.053635837568564, 1.053196121224978, 1.053032848771688, 1.053067868627288, 1.053182623885302, 1.053253336090298, 1.053665028109347, 1.053923372888378, 1.054074097007932, 1.054125496092429, 1.053895651205447, 1.053148829714424, 1.052679608454633, 1.052928319234431, 1.053908177329717, 1.055046653827171)
$colN = @(1.044718192082718, 1.045199529311993, 1.045510812323581, 1.045641632754068, 1.045663595526493, 1.045512560522839, 1.044880897631178, 1.04376654477692, 1.04302285252536, 1.042700653180671, 1.042580948806173, 1.042606626921506, 1.042690758886677, 1.042742592064962, 1.043044232235842, 1.043233396957714, 1.043343716507076, 1.04338132969687, 1.043213103126355, 1.042665984805621, 1.042321844645659, 1.042504293212891, 1.043222273087979, 1.044054776254914)

$columns = @{
    2 = $colB
    3 = $colC
    4 = $colD
    5 = $colE
    6 = $colF
    9 = $colI
    10 = $colJ
    11 = $colK
    12 = $colL
    13 = $colM
    14 = $colN
}

foreach ($colIndex in $columns.Keys) {
    $values = $columns[$colIndex]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = 2 + $i
        $ws.Cells.Item($row, $colIndex).Value = $values[$i]
    }
}
